$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ16127352"
$wb.Worksheets.Item(2).Name = "summ16251100"
$wb.Worksheets.Item(3).Name = "summ16349078"
$wb.Worksheets.Item(4).Name = "summ16458526"
$wb.Worksheets.Item(5).Name = "summ16568431"
$wb.Worksheets.Item(6).Name = "summ16736136"
$wb.Worksheets.Item(7).Name = "summ16853071"
$wb.Worksheets.Item(8).Name = "summ16995990"
$wb.Worksheets.Item(9).Name = "summ17109174"
$wb.Worksheets.Item(10).Name = "summ17217312"
$wb.Worksheets.Item(11).Name = "summ17323151"
$wb.Worksheets.Item(12).Name = "summ17435692"
$wb.Worksheets.Item(13).Name = "summ17540569"
$wb.Worksheets.Item(14).Name = "summ17758395"
$wb.Worksheets.Item(15).Name = "summ17885426"
$wb.Worksheets.Item(16).Name = "summ17990164"
$wb.Worksheets.Item(17).Name = "summ18093397"
$wb.Worksheets.Item(18).Name = "summ18205647"
$wb.Worksheets.Item(19).Name = "summ18315391"
$wb.Worksheets.Item(20).Name = "summ18420033"
$wb.Worksheets.Item(21).Name = "summ18530306"
$wb.Worksheets.Item(22).Name = "summ18636075"
$wb.Worksheets.Item(23).Name = "summ18739933"
$wb.Worksheets.Item(24).Name = "summ18849937"
$wb.Worksheets.Item(25).Name = "summ18954739"
$wb.Worksheets.Item(26).Name = "summ19065198"
$wb.Worksheets.Item(27).Name = "summ19169958"
$wb.Worksheets.Item(28).Name = "summ19280226"
$wb.Worksheets.Item(29).Name = "summ19386715"
$wb.Worksheets.Item(30).Name = "summ19499991"
$wb.Worksheets.Item(31).Name = "summ19609844"
$wb.Worksheets.Item(32).Name = "summ19719719"
$wb.Worksheets.Item(33).Name = "summ19835162"
$wb.Worksheets.Item(34).Name = "summ19939792"
$wb.Worksheets.Item(35).Name = "summ20050047"
$wb.Worksheets.Item(36).Name = "summ20156289"
$wb.Worksheets.Item(37).Name = "summ20275044"
$wb.Worksheets.Item(38).Name = "summ20395994"
$wb.Worksheets.Item(39).Name = "summ20562767"
$wb.Worksheets.Item(40).Name = "summ20704202"
$wb.Worksheets.Item(41).Name = "summ20849910"
$wb.Worksheets.Item(42).Name = "summ21004822"
$wb.Worksheets.Item(43).Name = "summ21145165"
$wb.Worksheets.Item(44).Name = "summ21304757"
$wb.Worksheets.Item(45).Name = "summ21459929"
$wb.Worksheets.Item(46).Name = "summ21604835"
$wb.Worksheets.Item(47).Name = "summ21745237"
$wb.Worksheets.Item(48).Name = "summ21900074"
$wb.Worksheets.Item(49).Name = "summ22060099"
$wb.Worksheets.Item(50).Name = "summ22225104"
